# Apply updated cryptocurrency market data to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers ("601.02", "0.141", ...) need to be
# forced to Text format first, otherwise Excel will auto-convert them to numeric values
# and lose the original text formatting (the source data are price strings, not numbers).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Now write the updated values.
$ws.Range('D2').Value = '66.641.65'
$ws.Range('E2').Value = '  -4.12%  '
$ws.Range('D3').Value = '3.453.45'
$ws.Range('E3').Value = '  -4.34%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '601.02'
$ws.Range('E5').Value = '  -4.53%  '
$ws.Range('D6').Value = '146.95'
$ws.Range('E6').Value = '  -7.50%  '
$ws.Range('D7').Value = '3.450.84'
$ws.Range('E7').Value = '  -4.40%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -2.53%  '
$ws.Range('D10').Value = '0.141'
$ws.Range('E10').Value = '  -5.48%  '
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('D12').Value = '0.421'
$ws.Range('E12').Value = '  -4.56%  '
$ws.Range('E13').Value = '  -7.51%  '
$ws.Range('D14').Value = '31.53'
$ws.Range('E14').Value = '  -5.95%  '
$ws.Range('D15').Value = '4.040.97'
$ws.Range('E15').Value = '  -4.22%  '
$ws.Range('D16').Value = '3.447.30'
$ws.Range('E16').Value = '  -4.36%  '
$ws.Range('D17').Value = '66.764.12'
$ws.Range('E17').Value = '  -3.98%  '
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').Value = '6.39'
$ws.Range('E19').Value = '  -4.63%  '
$ws.Range('E20').Value = '  -5.40%  '
$ws.Range('D21').Value = '9.95'
$ws.Range('E21').Value = '  -3.35%  '
$ws.Range('D22').Value = '437.59'
$ws.Range('E22').Value = '  -5.36%  '
$ws.Range('D23').Value = '0.606'
$ws.Range('E23').Value = '  -6.29%  '
$ws.Range('D24').Value = '78.00'
$ws.Range('E24').Value = '  -0.87%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').Value = '3.591.65'
$ws.Range('E26').Value = '  -4.28%  '
$ws.Range('D27').Value = '0.0000119'
$ws.Range('D28').Value = '9.82'
$ws.Range('E28').Value = '  -8.48%  '
$ws.Range('E29').Value = '  -9.85%  '
$ws.Range('E30').Value = '  -6.44%  '
$ws.Range('D31').Value = '1.59'
$ws.Range('E31').Value = '  -7.29%  '
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('E33').Value = '  -5.23%  '
$ws.Range('D34').Value = '25.32'
$ws.Range('E34').Value = '  -4.64%  '
$ws.Range('D35').Value = '6.08'
$ws.Range('E35').Value = '  -7.41%  '
$ws.Range('D36').Value = '3.448.90'
$ws.Range('E36').Value = '  -4.39%  '
$ws.Range('E37').Value = '  -7.92%  '
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('E39').Value = '  -7.36%  '
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('D41').Value = '173.13'
$ws.Range('E41').Value = '  -3.60%  '
$ws.Range('D42').Value = '2.16'
$ws.Range('E42').Value = '  -10.57%  '
$ws.Range('D43').Value = '0.0883'
$ws.Range('E43').Value = '  -4.82%  '
$ws.Range('D44').Value = '5.35'
$ws.Range('E44').Value = '  -5.44%  '
$ws.Range('E45').Value = '  -3.70%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '46.08'
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '28.66'
$ws.Range('E47').Value = '  -9.90%  '
$ws.Range('E48').Value = '  -11.88%  '
$ws.Range('D49').Value = '7.45'
$ws.Range('E49').Value = '  -4.81%  '
$ws.Range('D50').Value = '2.44'
$ws.Range('E50').Value = '  -11.22%  '
$ws.Range('D51').Value = '0.983'
$ws.Range('E51').Value = '  -5.21%  '
